$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column R ("Average Rent Prices")
$ws.Range("R1").Value = "Average Rent Prices"

# Row-wise average rent price (AVERAGE(C:Q)) as static values
$ws.Range("R2").Value = 191.53333333333333
$ws.Range("R3").Value = 189.6
$ws.Range("R4").Value = 187.66666666666666
$ws.Range("R5").Value = 190.4
$ws.Range("R6").Value = 189.06666666666666
$ws.Range("R7").Value = 185.93333333333334
$ws.Range("R8").Value = 191.8
$ws.Range("R9").Value = 196.46666666666667
$ws.Range("R10").Value = 202.93333333333334
$ws.Range("R11").Value = 192.33333333333334
$ws.Range("R12").Value = 192.13333333333333
$ws.Range("R13").Value = 183.86666666666667
$ws.Range("R14").Value = 180.53333333333333
$ws.Range("R15").Value = 170.2
$ws.Range("R16").Value = 169.93333333333334
$ws.Range("R17").Value = 166.4
$ws.Range("R18").Value = 158.66666666666666
$ws.Range("R19").Value = 151.80000000000001
$ws.Range("R20").Value = 150.33333333333334
$ws.Range("R21").Value = 148.06666666666666
$ws.Range("R22").Value = 156.93333333333334
$ws.Range("R23").Value = 159.4
$ws.Range("R24").Value = 168
$ws.Range("R25").Value = 168.93333333333334
$ws.Range("R26").Value = 177.4
$ws.Range("R27").Value = 171.73333333333332
$ws.Range("R28").Value = 188.13333333333333
$ws.Range("R29").Value = 193.06666666666666
$ws.Range("R30").Value = 189.06666666666666
$ws.Range("R31").Value = 191.73333333333332
$ws.Range("R32").Value = 203.8
$ws.Range("R33").Value = 199.73333333333332
$ws.Range("R34").Value = 207.46666666666667
$ws.Range("R35").Value = 212.26666666666668
$ws.Range("R36").Value = 229.46666666666667
$ws.Range("R37").Value = 236.8
$ws.Range("R38").Value = 253.86666666666667
$ws.Range("R39").Value = 261.53333333333336
$ws.Range("R40").Value = 288.2
$ws.Range("R41").Value = 237.8
$ws.Range("R42").Value = 212.2
$ws.Range("R43").Value = 201.93333333333334
$ws.Range("R44").Value = 225.6
$ws.Range("R45").Value = 231.26666666666668
$ws.Range("R46").Value = 239.06666666666666
$ws.Range("R47").Value = 255.6
$ws.Range("R48").Value = 263.60000000000002
$ws.Range("R49").Value = 267.86666666666667
$ws.Range("R50").Value = 270.86666666666667
$ws.Range("R51").Value = 278.2
$ws.Range("R52").Value = 299.86666666666667
$ws.Range("R53").Value = 282.2
$ws.Range("R54").Value = 278.86666666666667
$ws.Range("R55").Value = 286.60000000000002
$ws.Range("R56").Value = 279.46666666666664
$ws.Range("R57").Value = 293.2
$ws.Range("R58").Value = 286.46666666666664
$ws.Range("R59").Value = 286.39999999999998
$ws.Range("R60").Value = 317.66666666666669
$ws.Range("R61").Value = 315.93333333333334
$ws.Range("R62").Value = 307.53333333333336
$ws.Range("R63").Value = 312.93333333333334
$ws.Range("R64").Value = 329.53333333333336
$ws.Range("R65").Value = 326.13333333333333
$ws.Range("R66").Value = 334.33333333333331
$ws.Range("R67").Value = 341
$ws.Range("R68").Value = 352
$ws.Range("R69").Value = 333.06666666666666
$ws.Range("R70").Value = 322
$ws.Range("R71").Value = 317.46666666666664
$ws.Range("R72").Value = 331.13333333333333
$ws.Range("R73").Value = 328.4
$ws.Range("R74").Value = 338.86666666666667
$ws.Range("R75").Value = 342.53333333333336
$ws.Range("R76").Value = 354.86666666666667
$ws.Range("R77").Value = 356.26666666666665
$ws.Range("R78").Value = 362.8
$ws.Range("R79").Value = 367.4
$ws.Range("R80").Value = 372.6
$ws.Range("R81").Value = 370.33333333333331
$ws.Range("R82").Value = 367.4
$ws.Range("R83").Value = 365.26666666666665
$ws.Range("R84").Value = 376.93333333333334
$ws.Range("R85").Value = 363.73333333333335

